# Update sodan.xlsx ("相談件数" sheet) with the day's new figures and shift
# the footnote row down to make room for a blank spacer row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# --- 1. Fill in row 89 (2020-05-02 / serial 43944) with the day's counts ---
$ws.Range("B89").Value = 561
$ws.Range("C89").Value = 28381
$ws.Range("D89").Value = 159
$ws.Range("E89").Value = 6354

# --- 2. Insert a new blank row at 90, pushing the footnote row to 91 ---
$ws.Rows.Item(90).Insert()

# --- 3. Update the view: scroll/select so D90 is the active cell ---
$ws.Application.ActiveWindow.ScrollRow = 76
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("D90").Select()

# --- 4. Extend the sheet's print area by one row ($E$92 -> $E$93) ---
$wb.Names.Item(1).RefersTo = '=相談件数!$A$1:$E$93'

Write-Output "edit applied"
